# Update crypto price/volume data to the latest scraped snapshot.
# (Updated cryptos list on Sat Sep 14 15:11:38 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '59.947.48'
$ws.Range('E2').Value = '  +2.53%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.424.20'
$ws.Range('E3').Value = '  +2.29%  '

# Row 5: BNB
$ws.Range('D5').Value = '''552.41'
$ws.Range('E5').Value = '  +0.48%  '

# Row 6: Solana
$ws.Range('D6').Value = '''137.97'
$ws.Range('E6').Value = '  +2.94%  '

# Row 7: USDC
$ws.Range('E7').Value = '  -0.01%  '

# Row 8: XRP
$ws.Range('D8').Value = '''0.585'
$ws.Range('E8').Value = '  +2.54%  '

# Row 9: Dogecoin
$ws.Range('E9').Value = '  -1.03%  '

# Row 10: Toncoin
$ws.Range('E10').Value = '  +0.05%  '

# Row 11: Cardano
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').Value = '''0.356'
$ws.Range('E11').Value = '  -0.17%  '

# Row 12: TRON
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.148'
$ws.Range('E12').Value = '  -1.98%  '

# Row 13: Avalanche
$ws.Range('D13').Value = '''25.31'
$ws.Range('E13').Value = '  +4.12%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '2.855.31'
$ws.Range('E14').Value = '  +2.37%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '59.864.62'
$ws.Range('E15').Value = '  +2.53%  '

# Row 16: ShibaInu
$ws.Range('E16').Value = '  +0.70%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.397.52'
$ws.Range('E17').Value = '  +1.10%  '

# Row 18: Chainlink
$ws.Range('D18').Value = '''11.35'
$ws.Range('E18').Value = '  +2.04%  '

# Row 19: Polkadot
$ws.Range('E19').Value = '  +1.35%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = '''331.08'
$ws.Range('E20').Value = '  -0.55%  '

# Row 21: Uniswap
$ws.Range('E21').Value = '  -5.14%  '

# Row 22: Dai
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '  -0.16%  '

# Row 23: Litecoin
$ws.Range('D23').Value = '''66.28'
$ws.Range('E23').Value = '  +3.53%  '

# Row 24: Kaspa
$ws.Range('D24').Value = '''0.171'
$ws.Range('E24').Value = '  +1.26%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').Value = '''8.76'
$ws.Range('E25').Value = '  +5.00%  '

# Row 26: Binance-PegBSC-USD
$ws.Range('E26').Value = '  +0.31%  '

# Row 27: Fetch.AI
$ws.Range('E27').Value = '  +3.89%  '

# Row 28: PEPE
$ws.Range('D28').Value = '0.0₃0780'
$ws.Range('E28').Value = '  +4.25%  '

# Row 29: PancakeSwap
$ws.Range('E29').Value = '  +0.01%  '

# Row 30: Monero
$ws.Range('D30').Value = '''168.78'
$ws.Range('E30').Value = '  -0.90%  '

# Row 31: Aptos
$ws.Range('D31').Value = '''6.13'
$ws.Range('E31').Value = '  -0.85%  '

# Row 32: EthereumClassic
$ws.Range('E32').Value = '  +1.10%  '

# Row 33: SuiNetwork
$ws.Range('E33').Value = '  +1.34%  '

# Row 34: USDe
$ws.Range('E34').Value = '  +0.01%  '

# Row 35: ImmutableX
$ws.Range('E35').Value = '  +3.85%  '

# Row 36: NEARProtocol
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '''4.23'
$ws.Range('E36').Value = '  +0.95%  '

# Row 37: FirstDigitalUSD
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  +0.05%  '

# Row 38: Stacks
$ws.Range('E38').Value = '  +0.05%  '

# Row 39: OKB
$ws.Range('D39').Value = '''39.60'
$ws.Range('E39').Value = '  -2.14%  '

# Row 40: PolygonEcosystemToken
$ws.Range('E40').Value = '  -3.43%  '

# Row 41: Bittensor
$ws.Range('D41').Value = '''313.80'
$ws.Range('E41').Value = '  +8.08%  '

# Row 42: Filecoin
$ws.Range('E42').Value = '  -0.70%  '

# Row 43: Aave
$ws.Range('D43').Value = '''139.48'
$ws.Range('E43').Value = '  -1.17%  '

# Row 44: Stellar
$ws.Range('D44').Value = '''0.0969'
$ws.Range('E44').Value = '  +0.88%  '

# Row 45: Hedera
$ws.Range('D45').Value = '''0.0520'
$ws.Range('E45').Value = '  +0.66%  '

# Row 46: InjectiveProtocol
$ws.Range('D46').Value = '''19.55'
$ws.Range('E46').Value = '  +4.37%  '

# Row 47: Mantle
$ws.Range('E47').Value = '  +2.03%  '

# Row 48: VeChain
$ws.Range('E48').Value = '  +0.50%  '

# Row 49: Polygon
$ws.Range('D49').Value = '''0.392'
$ws.Range('E49').Value = '  -6.01%  '

# Row 50: EnergySwap
$ws.Range('D50').Value = '''17.63'
$ws.Range('E50').Value = '  +0.47%  '

# Row 51: WhiteBITCoin
$ws.Range('D51').Value = '''11.05'
$ws.Range('E51').Value = '  +0.17%  '
